$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.797.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.23%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.571.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.84%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'302.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'96.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.44%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.41%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.71%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.99%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.04%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.48%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +6.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.543.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.85%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.882"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.63%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'42.843.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.06%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'0.0₃0996"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.99%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").Value = "'12.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.90%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.83%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'72.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.33%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'253.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.80%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -5.38%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'28.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.21%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.10%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.67%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'37.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.84%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -5.84%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.10%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'154.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.62%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.14%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.56%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0801"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.52%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'18.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.64%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -3.07%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'22.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.30%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +29.56%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.88%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.52%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0310"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.30%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.082.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.95%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'9.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.11%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'85.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.92%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'75.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +9.10%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.25%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.821.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.74%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.191"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.85%  "
$ws.Range("E51").Style = "Normal"
